$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add week 44 (AU1) and week 45 (AV1) columns ---
# Force text storage (matches existing header cells D1:AT1, which are text "1".."43")
# while keeping the same bold/centered "General" style used across the header row.
$ws.Range("AU1").NumberFormat = "@"
$ws.Range("AU1").Value = "44"
$ws.Range("AV1").NumberFormat = "@"
$ws.Range("AV1").Value = "45"
$ws.Range("AT1").Copy() | Out-Null
$ws.Range("AU1:AV1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Company name correction (row 52) ---
$ws.Range("C52").Value = "COOMEVA EXPERIENCIA MEDICA SAS"

# --- New weekly data for week 44 (AU) and week 45 (AV), plus a few corrected/backfilled cells ---
$ws.Range("AU2").Value = 42
$ws.Range("AV2").Value = 44
$ws.Range("AV3").Value = 118
$ws.Range("AU5").Value = 5
$ws.Range("AV5").Value = 4
$ws.Range("AU6").Value = 54
$ws.Range("AV6").Value = 75
$ws.Range("AU7").Value = 15
$ws.Range("AV7").Value = 37
$ws.Range("AU8").Value = 11
$ws.Range("AV8").Value = 10
$ws.Range("AU9").Value = 1
$ws.Range("AV9").Value = 3
$ws.Range("AU10").Value = 2
$ws.Range("AU12").Value = 1
$ws.Range("AV12").Value = 5
$ws.Range("AV13").Value = 2
$ws.Range("AU14").Value = 2
$ws.Range("AV14").Value = 1
$ws.Range("AU15").Value = 2
$ws.Range("AV15").Value = 1
$ws.Range("AV16").Value = 1
$ws.Range("AV17").Value = 1
$ws.Range("AV22").Value = 1
$ws.Range("AU23").Value = 3
$ws.Range("AV23").Value = 4
$ws.Range("AU24").Value = 2
$ws.Range("AU25").Value = 28
$ws.Range("AV25").Value = 37
$ws.Range("AV26").Value = 2
$ws.Range("AU28").Value = 115
$ws.Range("AV28").Value = 140
$ws.Range("AU29").Value = 0
$ws.Range("AV29").Value = 0
$ws.Range("AU30").Value = 70
$ws.Range("AV30").Value = 52
$ws.Range("AU31").Value = 3
$ws.Range("AV31").Value = 0
$ws.Range("AU34").Value = 0
$ws.Range("AU35").Value = 33
$ws.Range("AV35").Value = 40
$ws.Range("AU36").Value = 1
$ws.Range("AV36").Value = 1
$ws.Range("E37").Value = 8
$ws.Range("F37").Value = 10
$ws.Range("G37").Value = 14
$ws.Range("H37").Value = 17
$ws.Range("I37").Value = 14
$ws.Range("AU37").Value = 6
$ws.Range("AV37").Value = 11
$ws.Range("AU38").Value = 95
$ws.Range("AV38").Value = 150
$ws.Range("AU40").Value = 1
$ws.Range("AT41").Value = 12
$ws.Range("AU41").Value = 5
$ws.Range("AV41").Value = 16
$ws.Range("AU42").Value = 30
$ws.Range("AV42").Value = 28
$ws.Range("AU43").Value = 105
$ws.Range("AV43").Value = 22
$ws.Range("AU44").Value = 373
$ws.Range("AU45").Value = 75
$ws.Range("AV45").Value = 48
$ws.Range("AU46").Value = 87
$ws.Range("AV46").Value = 103
$ws.Range("AU47").Value = 2
$ws.Range("AV47").Value = 3
$ws.Range("AU48").Value = 68
$ws.Range("AV48").Value = 106
$ws.Range("AU49").Value = 2
$ws.Range("AV49").Value = 4
$ws.Range("AU50").Value = 0
$ws.Range("AV50").Value = 0
$ws.Range("AU51").Value = 3
$ws.Range("AU53").Value = 3
$ws.Range("AV53").Value = 17
$ws.Range("AU54").Value = 0
$ws.Range("AV54").Value = 2
$ws.Range("AU55").Value = 0
$ws.Range("AV55").Value = 1
$ws.Range("AU56").Value = 6
$ws.Range("AV56").Value = 7
$ws.Range("AU57").Value = 18
$ws.Range("AV57").Value = 108
$ws.Range("AU58").Value = 24
$ws.Range("AV58").Value = 12
